$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 1100.1487
$ws.Range("I15").Value = 1100.1487
$ws.Range("K15").Value = 3300.4461
$ws.Range("M15").Value = -3131.4461

# ALC row 19
$ws.Range("H19").Value = 3745.0952
$ws.Range("I19").Value = 5044.6665
$ws.Range("J19").Value = 2012.3334
$ws.Range("K19").Value = 5044.6665
$ws.Range("L19").Value = 2012.3334
$ws.Range("M19").Value = -4869.6665
$ws.Range("N19").Value = -2362.3334

# ALC row 98
$ws.Range("H98").Value = 1143.3871
$ws.Range("I98").Value = 1164.8334
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 1164.8334
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 333.1666
$ws.Range("N98").Value = -3496

# ALC row 122
$ws.Range("H122").Value = 1143.3871
$ws.Range("I122").Value = 1164.8334
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 3494.5002
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -1044.5002
$ws.Range("N122").Value = -6400

# ALC row 137
$ws.Range("H137").Value = 1287.0555
$ws.Range("I137").Value = 1156.2142
$ws.Range("J137").Value = 1745
$ws.Range("K137").Value = 3468.6426
$ws.Range("L137").Value = 5235
$ws.Range("M137").Value = -918.6425999999997
$ws.Range("N137").Value = -10335

# ALC row 138
$ws.Range("H138").Value = 1566.725
$ws.Range("I138").Value = 991.1539
$ws.Range("J138").Value = 2635.6428
$ws.Range("K138").Value = 2973.4617
$ws.Range("L138").Value = 7906.928400000001
$ws.Range("M138").Value = 2166.5383
$ws.Range("N138").Value = -18186.9284

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 18817.662
$ws.Range("I32").Value = 20203.422
$ws.Range("J32").Value = 3020
$ws.Range("K32").Value = 20203.422
$ws.Range("L32").Value = 3020
$ws.Range("M32").Value = -19916.422
$ws.Range("N32").Value = -3594

# ARM row 107
$ws.Range("H107").Value = 57575
$ws.Range("J107").Value = 57575
$ws.Range("L107").Value = 57575
$ws.Range("N107").Value = -65255

# ARM row 109
$ws.Range("H109").Value = 22500
$ws.Range("J109").Value = 22500
$ws.Range("L109").Value = 22500
$ws.Range("N109").Value = -25274

$ws = $wb.Worksheets.Item("BSM")
# BSM row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

# BSM row 107
$ws.Range("H107").Value = 1636.7142
$ws.Range("I107").Value = 903.6667
$ws.Range("J107").Value = 1836.6364
$ws.Range("K107").Value = 903.6667
$ws.Range("L107").Value = 1836.6364
$ws.Range("M107").Value = 1016.3333
$ws.Range("N107").Value = -5676.6364

$ws = $wb.Worksheets.Item("CUL")
# CUL row 58
$ws.Range("H58").Value = 2755.75
$ws.Range("I58").Value = 2005
$ws.Range("J58").Value = 3006
$ws.Range("K58").Value = 6015
$ws.Range("L58").Value = 9018
$ws.Range("M58").Value = -5887
$ws.Range("N58").Value = -9274

# CUL row 69
$ws.Range("H69").Value = 1830.6897

# CUL row 72
$ws.Range("H72").Value = 1830.6897

# CUL row 131
$ws.Range("H131").Value = 3240.2144
$ws.Range("I131").Value = 14742.857
$ws.Range("J131").Value = 939.6857
$ws.Range("K131").Value = 44228.571
$ws.Range("L131").Value = 2819.0571
$ws.Range("M131").Value = -39188.571
$ws.Range("N131").Value = -12899.0571

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Range("H122").Value = 9650
$ws.Range("I122").Value = 29600
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 88800
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -86350
$ws.Range("N122").Value = -13900

# GSM row 132
$ws.Range("H132").Value = 43106.32
$ws.Range("I132").Value = 64834.188
$ws.Range("K132").Value = 194502.564
$ws.Range("M132").Value = -191972.564

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 2534.4211
$ws.Range("I7").Value = 1501.3334
$ws.Range("J7").Value = 2728.125
$ws.Range("K7").Value = 1501.3334
$ws.Range("L7").Value = 2728.125
$ws.Range("M7").Value = -1389.3334
$ws.Range("N7").Value = -2952.125

# LTW row 76
$ws.Range("H76").Value = 14216
$ws.Range("J76").Value = 14144
$ws.Range("L76").Value = 14144
$ws.Range("N76").Value = -14820

# LTW row 79
$ws.Range("H79").Value = 14216
$ws.Range("J79").Value = 14144
$ws.Range("L79").Value = 14144
$ws.Range("N79").Value = -16484

# LTW row 93
$ws.Range("H93").Value = 1802955.2
$ws.Range("I93").Value = 3862164.5
$ws.Range("J93").Value = 1147.25
$ws.Range("K93").Value = 3862164.5
$ws.Range("L93").Value = 1147.25
$ws.Range("M93").Value = -3860916.5
$ws.Range("N93").Value = -3643.25

# LTW row 122
$ws.Range("H122").Value = 3046.1177
$ws.Range("I122").Value = 3291.2727
$ws.Range("J122").Value = 2596.6667
$ws.Range("K122").Value = 9873.8181
$ws.Range("L122").Value = 7790.000100000001
$ws.Range("M122").Value = -7423.8181
$ws.Range("N122").Value = -12690.0001

# LTW row 126
$ws.Range("H126").Value = 2534.4211
$ws.Range("I126").Value = 1501.3334
$ws.Range("J126").Value = 2728.125
$ws.Range("K126").Value = 4504.0002
$ws.Range("L126").Value = 8184.375
$ws.Range("M126").Value = -2034.0002
$ws.Range("N126").Value = -13124.375

$ws = $wb.Worksheets.Item("WVR")
# WVR row 62
$ws.Range("H62").Value = 4513.3335
$ws.Range("I62").Value = 3980
$ws.Range("J62").Value = 4620
$ws.Range("K62").Value = 3980
$ws.Range("L62").Value = 4620
$ws.Range("M62").Value = -3356
$ws.Range("N62").Value = -5868

# WVR row 63
$ws.Range("H63").Value = 14900
$ws.Range("I63").Value = 14900
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 14900
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("M63").Value = -14276

# WVR row 65
$ws.Range("H65").Value = 4513.3335
$ws.Range("I65").Value = 3980
$ws.Range("J65").Value = 4620
$ws.Range("K65").Value = 19900
$ws.Range("L65").Value = 23100
$ws.Range("M65").Value = -16780
$ws.Range("N65").Value = -29340

# WVR row 66
$ws.Range("H66").Value = 14900
$ws.Range("I66").Value = 14900
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 44700
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("M66").Value = -41580

# WVR row 76
$ws.Range("H76").Value = 18724.334
$ws.Range("J76").Value = 18586.5
$ws.Range("L76").Value = 18586.5
$ws.Range("N76").Value = -19216.5

# WVR row 79
$ws.Range("H79").Value = 18724.334
$ws.Range("J79").Value = 18586.5
$ws.Range("L79").Value = 18586.5
$ws.Range("N79").Value = -20770.5

# WVR row 122
$ws.Range("H122").Value = 1871.7142
$ws.Range("I122").Value = 2378.2222
$ws.Range("J122").Value = 960
$ws.Range("K122").Value = 7134.6666
$ws.Range("L122").Value = 2880
$ws.Range("M122").Value = -4684.6666
$ws.Range("N122").Value = -7780
